$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 11.90597433333333
$ws.Range("H2").Value = 35.717923
$ws.Range("I2").Value = 0.008895149679642379
$ws.Range("J2").Value = 0.008895149679642379
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 7.321929333333333
$ws.Range("N2").Value = 21.965788
$ws.Range("Q2").Value = 87.17470271314711
$ws.Range("R2").Value = 784.572324418324
$ws.Range("S2").Value = 0.008895149679642379
$ws.Range("T2").Value = 0.008895149679642379

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1205.102620666667
$ws.Range("H3").Value = 3615.307862
$ws.Range("I3").Value = 0.9003520325209805
$ws.Range("J3").Value = 0.9003520325209804
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.321929333333333
$ws.Range("N3").Value = 21.965788
$ws.Range("Q3").Value = 8823.676227936141
$ws.Range("R3").Value = 79413.08605142526
$ws.Range("S3").Value = 0.9003520325209805
$ws.Range("T3").Value = 0.9003520325209804

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 121.4707743333333
$ws.Range("H4").Value = 364.412323
$ws.Range("I4").Value = 0.09075281779937723
$ws.Range("J4").Value = 0.09075281779937722
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 7.321929333333333
$ws.Range("N4").Value = 21.965788
$ws.Range("Q4").Value = 889.4004257339471
$ws.Range("R4").Value = 8004.603831605525
$ws.Range("S4").Value = 0.09075281779937723
$ws.Range("T4").Value = 0.09075281779937722
